$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9773569703102112
$ws.Range("B1").Value = 1.852310061454773
$ws.Range("C1").Value = 3.301288366317749
$ws.Range("D1").Value = 3.928803205490112
$ws.Range("E1").Value = 0.4214289784431458
